# Auto-generated edit script: updates cached market-price values
# across the 8 job sheets (scheduled-runner style refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 607.907
$ws.Range("J17").Value = 591.5294
$ws.Range("L17").Value = 1774.5882
$ws.Range("N17").Value = -2110.5882

$ws.Range("H28").Value = 30303542
$ws.Range("I28").Value = 37037324
$ws.Range("K28").Value = 37037324
$ws.Range("M28").Value = -37036839

$ws.Range("H98").Value = 1175.5946
$ws.Range("I98").Value = 1066.9667
$ws.Range("K98").Value = 1066.9667
$ws.Range("M98").Value = 431.0333000000001

$ws.Range("H106").Value = 15878627
$ws.Range("I106").Value = 16672058
$ws.Range("K106").Value = 16672058
$ws.Range("M106").Value = -16671427

$ws.Range("H116").Value = 10003.333
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 10003.333
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 10003.333
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = -16887.333

$ws.Range("H122").Value = 1175.5946
$ws.Range("I122").Value = 1066.9667
$ws.Range("K122").Value = 3200.9001
$ws.Range("M122").Value = -750.9000999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 44049332
$ws.Range("I2").Value = 48053588
$ws.Range("K2").Value = 48053588
$ws.Range("M2").Value = -48053475

$ws.Range("H45").Value = 14088.637
$ws.Range("I45").Value = 27481
$ws.Range("K45").Value = 27481
$ws.Range("M45").Value = -27104

$ws.Range("H74").Value = 17609.223
$ws.Range("I74").Value = 17609.223
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 17609.223
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -16735.223
$ws.Range("N74").Value = ""

$ws.Range("H77").Value = 17609.223
$ws.Range("I77").Value = 17609.223
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 88046.11500000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -83678.11500000001
$ws.Range("N77").Value = ""

$ws.Range("H116").Value = 44049332
$ws.Range("I116").Value = 48053588
$ws.Range("K116").Value = 48053588
$ws.Range("M116").Value = -48051294

$ws.Range("H122").Value = 2358.52
$ws.Range("I122").Value = 1776.4706
$ws.Range("K122").Value = 5329.4118
$ws.Range("M122").Value = -2879.4118

$ws.Range("H132").Value = 5542.2383
$ws.Range("I132").Value = 4764.2
$ws.Range("K132").Value = 14292.6
$ws.Range("M132").Value = -11762.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 44049332
$ws.Range("I3").Value = 48053588
$ws.Range("K3").Value = 48053588
$ws.Range("M3").Value = -48053474

$ws.Range("H94").Value = 2395.5
$ws.Range("I94").Value = 2243.375
$ws.Range("K94").Value = 2243.375
$ws.Range("M94").Value = -1792.375

$ws.Range("H105").Value = 2852.375
$ws.Range("I105").Value = 2909.2
$ws.Range("K105").Value = 2909.2
$ws.Range("M105").Value = -1162.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 77500
$ws.Range("I20").Value = 60000
$ws.Range("J20").Value = 80000
$ws.Range("K20").Value = 60000
$ws.Range("L20").Value = 80000
$ws.Range("M20").Value = -59764
$ws.Range("N20").Value = -80472

$ws.Range("H30").Value = 77500
$ws.Range("I30").Value = 60000
$ws.Range("J30").Value = 80000
$ws.Range("K30").Value = 60000
$ws.Range("L30").Value = 80000
$ws.Range("M30").Value = -59909
$ws.Range("N30").Value = -80182

$ws.Range("H31").Value = 1738.6364
$ws.Range("I31").Value = 1539.5264
$ws.Range("J31").Value = 2999.6667
$ws.Range("K31").Value = 1539.5264
$ws.Range("L31").Value = 2999.6667
$ws.Range("M31").Value = -1244.5264
$ws.Range("N31").Value = -3589.6667

$ws.Range("H34").Value = 1738.6364
$ws.Range("I34").Value = 1539.5264
$ws.Range("J34").Value = 2999.6667
$ws.Range("K34").Value = 1539.5264
$ws.Range("L34").Value = 2999.6667
$ws.Range("M34").Value = -1337.5264
$ws.Range("N34").Value = -3403.6667

$ws.Range("H59").Value = 95000
$ws.Range("I59").Value = 95000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 95000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -93855
$ws.Range("N59").Value = ""

$ws.Range("H99").Value = 9141.65
$ws.Range("I99").Value = 10778.444
$ws.Range("J99").Value = 8666.450999999999
$ws.Range("K99").Value = 10778.444
$ws.Range("L99").Value = 8666.450999999999
$ws.Range("M99").Value = -9280.444
$ws.Range("N99").Value = -11662.451

$ws.Range("H122").Value = 2851.7144
$ws.Range("I122").Value = 3009.25
$ws.Range("J122").Value = 2641.6667
$ws.Range("K122").Value = 9027.75
$ws.Range("L122").Value = 7925.000100000001
$ws.Range("M122").Value = -6577.75
$ws.Range("N122").Value = -12825.0001

$ws.Range("H126").Value = 9141.65
$ws.Range("I126").Value = 10778.444
$ws.Range("J126").Value = 8666.450999999999
$ws.Range("K126").Value = 32335.332
$ws.Range("L126").Value = 25999.353
$ws.Range("M126").Value = -29865.332
$ws.Range("N126").Value = -30939.353

$ws.Range("H128").Value = 77500
$ws.Range("I128").Value = 60000
$ws.Range("J128").Value = 80000
$ws.Range("K128").Value = 60000
$ws.Range("L128").Value = 80000
$ws.Range("M128").Value = -55020
$ws.Range("N128").Value = -89960

$ws.Range("H132").Value = 3365.4614
$ws.Range("I132").Value = 2761
$ws.Range("J132").Value = 7999.6665
$ws.Range("K132").Value = 8283
$ws.Range("L132").Value = 23998.9995
$ws.Range("M132").Value = -5753
$ws.Range("N132").Value = -29058.9995

$ws.Range("H134").Value = 5935.879
$ws.Range("I134").Value = 6289.4
$ws.Range("J134").Value = 4831.125
$ws.Range("K134").Value = 18868.2
$ws.Range("L134").Value = 14493.375
$ws.Range("M134").Value = -16333.2
$ws.Range("N134").Value = -19563.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1021.5
$ws.Range("I113").Value = 821.2857
$ws.Range("J113").Value = 1114.9333
$ws.Range("K113").Value = 2463.8571
$ws.Range("L113").Value = 3344.7999
$ws.Range("M113").Value = -293.8571000000002
$ws.Range("N113").Value = -7684.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5827
$ws.Range("I122").Value = 6427.857
$ws.Range("J122").Value = 4425
$ws.Range("K122").Value = 19283.571
$ws.Range("L122").Value = 13275
$ws.Range("M122").Value = -16833.571
$ws.Range("N122").Value = -18175

$ws.Range("H126").Value = 3887.8113
$ws.Range("I126").Value = 4549.5625
$ws.Range("K126").Value = 13648.6875
$ws.Range("M126").Value = -11178.6875

$ws.Range("H132").Value = 12617.333
$ws.Range("I132").Value = 16611.777
$ws.Range("J132").Value = 6625.6665
$ws.Range("K132").Value = 49835.33099999999
$ws.Range("L132").Value = 19876.9995
$ws.Range("M132").Value = -47305.33099999999
$ws.Range("N132").Value = -24936.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8931234
$ws.Range("I16").Value = 9617490
$ws.Range("J16").Value = 9899
$ws.Range("K16").Value = 9617490
$ws.Range("L16").Value = 9899
$ws.Range("M16").Value = -9617320
$ws.Range("N16").Value = -10239

$ws.Range("H22").Value = 979.8333
$ws.Range("I22").Value = 971.25
$ws.Range("J22").Value = 997
$ws.Range("K22").Value = 971.25
$ws.Range("L22").Value = 997
$ws.Range("M22").Value = -676.25
$ws.Range("N22").Value = -1587

$ws.Range("H27").Value = 979.8333
$ws.Range("I27").Value = 971.25
$ws.Range("J27").Value = 997
$ws.Range("K27").Value = 971.25
$ws.Range("L27").Value = 997
$ws.Range("M27").Value = -864.25
$ws.Range("N27").Value = -1211

$ws.Range("H61").Value = 6656.136
$ws.Range("I61").Value = 7890.4116
$ws.Range("J61").Value = 2459.6
$ws.Range("K61").Value = 7890.4116
$ws.Range("L61").Value = 2459.6
$ws.Range("M61").Value = -7688.4116
$ws.Range("N61").Value = -2863.6

$ws.Range("H113").Value = 6656.136
$ws.Range("I113").Value = 7890.4116
$ws.Range("J113").Value = 2459.6
$ws.Range("K113").Value = 7890.4116
$ws.Range("L113").Value = 2459.6
$ws.Range("M113").Value = -5720.4116
$ws.Range("N113").Value = -6799.6

$ws.Range("H122").Value = 3600
$ws.Range("I122").Value = 2400
$ws.Range("K122").Value = 7200
$ws.Range("M122").Value = -4750

$ws.Range("H134").Value = 96353.75
$ws.Range("J134").Value = 96353.75
$ws.Range("L134").Value = 96353.75
$ws.Range("N134").Value = -106493.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11198.553
$ws.Range("I122").Value = 8710.138000000001
$ws.Range("K122").Value = 26130.414
$ws.Range("M122").Value = -23680.414

$ws.Range("H133").Value = 74897.60000000001
$ws.Range("J133").Value = 74897.60000000001
$ws.Range("L133").Value = 74897.60000000001
$ws.Range("N133").Value = -85017.60000000001

$ws.Range("H136").Value = 2006.3684
$ws.Range("I136").Value = 2001.2285
$ws.Range("J136").Value = 2066.3333
$ws.Range("K136").Value = 6003.6855
$ws.Range("L136").Value = 6198.999899999999
$ws.Range("M136").Value = -3453.6855
$ws.Range("N136").Value = -11298.9999
